# GUI: Updated the statistics.
#
# A new test ("NestedFolders") was added to the tracked-tests table on
# Sheet1. This pushes the two rows that used to sit at rows 10-11 down
# to rows 11-12, and the new test's data is written into row 10.
#
# We deliberately do the row "insert" by copying values down one row at
# a time (rather than Rows.Insert, which would also shift every
# formatting-only row further down the sheet) so that only rows 10-12
# change, matching the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old row 11 (Overlay / 3 / 4 / Finished) moves down to row 12.
$ws.Range("A12").Value = "Overlay"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "Finished"

# Old row 10 (Non-CustodianObjects / 0 / 1 / Ready to Write) moves down to row 11.
$ws.Range("A11").Value = "Non-CustodianObjects"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Ready to Write"

# New row 10 holds the newly added test.
$ws.Range("A10").Value = "NestedFolders"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "Automated"

# Update the active selection to match where the user last clicked.
$ws.Range("E10").Select()
